# Scheduled market-data refresh: update computed price/profit values
# across the Leve profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 6686.5264
$ws.Range("J28").Value = 9931.75
$ws.Range("L28").Value = 9931.75
$ws.Range("N28").Value = -10901.75
$ws.Range("H42").Value = 648.25
$ws.Range("I42").Value = 726.7143
$ws.Range("J42").Value = 99
$ws.Range("K42").Value = 2180.1429
$ws.Range("L42").Value = 297
$ws.Range("M42").Value = -1950.1429
$ws.Range("N42").Value = -757
$ws.Range("H70").Value = 1984.0769
$ws.Range("I70").Value = 1179.4
$ws.Range("J70").Value = 2487
$ws.Range("K70").Value = 3538.2
$ws.Range("L70").Value = 7461
$ws.Range("M70").Value = -3268.2
$ws.Range("N70").Value = -8001
$ws.Range("H73").Value = 1984.0769
$ws.Range("I73").Value = 1179.4
$ws.Range("J73").Value = 2487
$ws.Range("K73").Value = 3538.2
$ws.Range("L73").Value = 7461
$ws.Range("M73").Value = -2602.2
$ws.Range("N73").Value = -9333
$ws.Range("H76").Value = 58828700
$ws.Range("I76").Value = 58828700
$ws.Range("K76").Value = 58828700
$ws.Range("M76").Value = -58828385
$ws.Range("H79").Value = 58828700
$ws.Range("I79").Value = 58828700
$ws.Range("K79").Value = 58828700
$ws.Range("M79").Value = -58827608
$ws.Range("H80").Value = 215.72728
$ws.Range("I80").Value = 284.42856
$ws.Range("J80").Value = 95.5
$ws.Range("K80").Value = 853.28568
$ws.Range("L80").Value = 286.5
$ws.Range("M80").Value = 144.71432
$ws.Range("N80").Value = -2282.5
$ws.Range("H82").Value = 4168.8
$ws.Range("I82").Value = 4168.8
$ws.Range("K82").Value = 12506.4
$ws.Range("M82").Value = -12100.4
$ws.Range("H83").Value = 215.72728
$ws.Range("I83").Value = 284.42856
$ws.Range("J83").Value = 95.5
$ws.Range("K83").Value = 2559.85704
$ws.Range("L83").Value = 859.5
$ws.Range("M83").Value = 2432.14296
$ws.Range("N83").Value = -10843.5
$ws.Range("H85").Value = 4168.8
$ws.Range("I85").Value = 4168.8
$ws.Range("K85").Value = 12506.4
$ws.Range("M85").Value = -11102.4
$ws.Range("H100").Value = 5875
$ws.Range("I100").Value = 5875
$ws.Range("K100").Value = 5875
$ws.Range("M100").Value = -5334
$ws.Range("H101").Value = 403.875
$ws.Range("I101").Value = 425.85715
$ws.Range("K101").Value = 1277.57145
$ws.Range("M101").Value = 344.4285500000001
$ws.Range("H132").Value = 253993.42
$ws.Range("J132").Value = 9680
$ws.Range("L132").Value = 29040
$ws.Range("N132").Value = -34100
$ws.Range("H135").Value = 4665.4
$ws.Range("I135").Value = 2215.8845
$ws.Range("J135").Value = 9214.5
$ws.Range("K135").Value = 19942.9605
$ws.Range("L135").Value = 82930.5
$ws.Range("M135").Value = -17407.9605
$ws.Range("N135").Value = -88000.5
$ws.Range("H138").Value = 3923.775
$ws.Range("I138").Value = 1511.2858
$ws.Range("K138").Value = 4533.857400000001
$ws.Range("M138").Value = 606.1425999999992

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1913.3334
$ws.Range("J97").Value = 1898
$ws.Range("L97").Value = 1898
$ws.Range("N97").Value = -2890
$ws.Range("H122").Value = 2568.6956
$ws.Range("I122").Value = 1890
$ws.Range("K122").Value = 5670
$ws.Range("M122").Value = -3220
$ws.Range("H132").Value = 3045824.2
$ws.Range("I132").Value = 10770691
$ws.Range("J132").Value = 148999.5
$ws.Range("K132").Value = 32312073
$ws.Range("L132").Value = 446998.5
$ws.Range("M132").Value = -32309543
$ws.Range("N132").Value = -452058.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 8344.227000000001
$ws.Range("I99").Value = 7930.8647
$ws.Range("K99").Value = 7930.8647
$ws.Range("M99").Value = -6432.8647
$ws.Range("H134").Value = 3346143.8
$ws.Range("J134").Value = 20749.5
$ws.Range("L134").Value = 62248.5
$ws.Range("N134").Value = -67318.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 17244812
$ws.Range("I16").Value = 45456364
$ws.Range("K16").Value = 45456364
$ws.Range("M16").Value = -45456077
$ws.Range("H62").Value = 5992.077
$ws.Range("I62").Value = 6329.7
$ws.Range("J62").Value = 4866.6665
$ws.Range("K62").Value = 6329.7
$ws.Range("L62").Value = 4866.6665
$ws.Range("M62").Value = -5705.7
$ws.Range("N62").Value = -6114.6665
$ws.Range("H65").Value = 5992.077
$ws.Range("I65").Value = 6329.7
$ws.Range("J65").Value = 4866.6665
$ws.Range("K65").Value = 31648.5
$ws.Range("L65").Value = 24333.3325
$ws.Range("M65").Value = -28528.5
$ws.Range("N65").Value = -30573.3325
$ws.Range("H105").Value = 62501280
$ws.Range("I105").Value = 90910460
$ws.Range("K105").Value = 90910460
$ws.Range("M105").Value = -90908713
$ws.Range("H113").Value = 17244812
$ws.Range("I113").Value = 45456364
$ws.Range("K113").Value = 45456364
$ws.Range("M113").Value = -45454194
$ws.Range("H122").Value = 2703.7368
$ws.Range("I122").Value = 2147.3333
$ws.Range("K122").Value = 6441.999899999999
$ws.Range("M122").Value = -3991.999899999999
$ws.Range("H132").Value = 14659.637
$ws.Range("I132").Value = 12419.444
$ws.Range("K132").Value = 37258.33199999999
$ws.Range("M132").Value = -34728.33199999999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 270
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H97").Value = 477.4
$ws.Range("I97").Value = 463.33334
$ws.Range("J97").Value = 498.5
$ws.Range("K97").Value = 1390.00002
$ws.Range("L97").Value = 1495.5
$ws.Range("M97").Value = -894.0000199999999
$ws.Range("N97").Value = -2487.5
$ws.Range("H121").Value = 75017.164
$ws.Range("I121").Value = 2684
$ws.Range("J121").Value = 121047.37
$ws.Range("K121").Value = 8052
$ws.Range("L121").Value = 363142.11
$ws.Range("M121").Value = -6742
$ws.Range("N121").Value = -365762.11

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 12999.714
$ws.Range("I113").Value = 6999
$ws.Range("K113").Value = 6999
$ws.Range("M113").Value = -4829
$ws.Range("H122").Value = 7867.2
$ws.Range("J122").Value = 8712.143
$ws.Range("L122").Value = 26136.429
$ws.Range("N122").Value = -31036.429
$ws.Range("H132").Value = 37042404
$ws.Range("I132").Value = 55561040
$ws.Range("K132").Value = 166683120
$ws.Range("M132").Value = -166680590

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6097.72
$ws.Range("I61").Value = 4812.4116
$ws.Range("K61").Value = 4812.4116
$ws.Range("M61").Value = -4610.4116
$ws.Range("H68").Value = 2384.6191
$ws.Range("I68").Value = 1879.75
$ws.Range("K68").Value = 1879.75
$ws.Range("M68").Value = -1130.75
$ws.Range("H71").Value = 2384.6191
$ws.Range("I71").Value = 1879.75
$ws.Range("K71").Value = 9398.75
$ws.Range("M71").Value = -5654.75
$ws.Range("H113").Value = 6097.72
$ws.Range("I113").Value = 4812.4116
$ws.Range("K113").Value = 4812.4116
$ws.Range("M113").Value = -2642.4116

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 11862.2
$ws.Range("I62").Value = 10649.5
$ws.Range("K62").Value = 10649.5
$ws.Range("M62").Value = -10025.5
$ws.Range("H65").Value = 11862.2
$ws.Range("I65").Value = 10649.5
$ws.Range("K65").Value = 53247.5
$ws.Range("M65").Value = -50127.5
$ws.Range("H81").Value = 1247.7059
$ws.Range("I81").Value = 1142.75
$ws.Range("J81").Value = 1499.6
$ws.Range("K81").Value = 2285.5
$ws.Range("L81").Value = 2999.2
$ws.Range("M81").Value = -1224.5
$ws.Range("N81").Value = -5121.2
$ws.Range("H84").Value = 1247.7059
$ws.Range("I84").Value = 1142.75
$ws.Range("J84").Value = 1499.6
$ws.Range("K84").Value = 11427.5
$ws.Range("L84").Value = 14996
$ws.Range("M84").Value = -6123.5
$ws.Range("N84").Value = -25604
$ws.Range("H113").Value = 7577567
$ws.Range("I113").Value = 11906539
$ws.Range("K113").Value = 35719617
$ws.Range("M113").Value = -35717447
$ws.Range("H122").Value = 6743.2163
$ws.Range("I122").Value = 3196.5833
$ws.Range("J122").Value = 13290.846
$ws.Range("K122").Value = 9589.749899999999
$ws.Range("L122").Value = 39872.538
$ws.Range("M122").Value = -7139.749899999999
$ws.Range("N122").Value = -44772.538
